# Weekly fruit/vegetable price update: insert a new weekly record for
# "Apio" (Vega Modelo de Temuco) ahead of the existing series, pushing the
# rest of the historical rows (old 169..201) down by one row (new 170..202).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 169 — shifts rows 169:201 down to 170:202 and
# extends the sheet's used range to A1:R202, matching the diff's dimension.
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A169").Value = 10
$ws.Range("B169").Value = "Vega Modelo de Temuco"
$ws.Range("C169").Value = "La Araucanía"
$ws.Range("D169").Value = 44504
$ws.Range("E169").Value = 9
$ws.Range("F169").Value = 100112017
$ws.Range("G169").Value = "Apio"
$ws.Range("H169").Value = "Americana (o)"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 220
$ws.Range("K169").Value = 8000
$ws.Range("L169").Value = 9000
$ws.Range("M169").Value = 8568
$ws.Range("N169").Value = "`$/docena de matas"
$ws.Range("O169").Value = "Provincia del Elquí"
$ws.Range("P169").Value = 1428
$ws.Range("Q169").Value = 6
$ws.Range("R169").Value = "Hortaliza"
